$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three data rows (2,3,4) get cyclically rotated:
#   new row 2 <- old row 4
#   new row 3 <- old row 2
#   new row 4 <- old row 3
# Columns A (Id), Q (Ost) and R (Nord) carry the only differing values
# between these rows, so read the "before" values first and then write
# them back in rotated order.

$a2 = $ws.Range("A2").Value2
$a3 = $ws.Range("A3").Value2
$a4 = $ws.Range("A4").Value2

$q2 = $ws.Range("Q2").Value2
$q3 = $ws.Range("Q3").Value2
$q4 = $ws.Range("Q4").Value2

$r2 = $ws.Range("R2").Value2
$r3 = $ws.Range("R3").Value2
$r4 = $ws.Range("R4").Value2

$ws.Range("A2").Value = $a4
$ws.Range("A3").Value = $a2
$ws.Range("A4").Value = $a3

$ws.Range("Q2").Value = $q4
$ws.Range("Q3").Value = $q2
$ws.Range("Q4").Value = $q3

$ws.Range("R2").Value = $r4
$ws.Range("R3").Value = $r2
$ws.Range("R4").Value = $r3

# Row 2 previously carried a handful of blank attribute cells (J,K,L,N,AF)
# that belonged to the record that has now moved to row 4 - those blank
# cells move along with the rest of the record, so clear them from row 2...
$ws.Range("J2").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("AF2").ClearContents()

# ...and (re)create them as blank cells on row 4, matching what row 2 had.
$ws.Range("J4").Style = "Normal"
$ws.Range("K4").Style = "Normal"
$ws.Range("L4").Style = "Normal"
$ws.Range("N4").Style = "Normal"
$ws.Range("AF4").Style = "Normal"
